$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 15: new labels with wrapped / centered formatting ---
$ws.Range("C15").Value = "Total gas (Mscf/day)"
$ws.Range("D15").Value = "Total oil production (bbl/day )"

$hdr = $ws.Range("C15:D15")
$hdr.Font.Name = "Times New Roman"
$hdr.Font.Size = 10
$hdr.Font.Color = 0
$hdr.Interior.Color = 16777215
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108
$hdr.WrapText = $true
$hdr.Borders.LineStyle = 1
$ws.Rows(15).RowHeight = 51

# --- Row 16: totals (gas + C1 fraction; oil totals) ---
$ws.Range("A16").Formula = "=A3+I3"
$ws.Range("C16").Formula = "=C3+K3"
$ws.Range("D16").Formula = "=D3+L3"

# --- Rows 17-25: shared-formula fill down from row 16 ---
$ws.Range("A17:A25").Formula = "=A4+I4"
$ws.Range("C17:D17").Formula = "=C4+K4"
$ws.Range("C18:D18").Formula = "=C5+K5"
$ws.Range("C19:D19").Formula = "=C6+K6"
$ws.Range("C20:D20").Formula = "=C7+K7"
$ws.Range("C21:D21").Formula = "=C8+K8"
$ws.Range("C22:D22").Formula = "=C9+K9"
$ws.Range("C23:D23").Formula = "=C10+K10"
$ws.Range("C24:D24").Formula = "=C11+K11"
$ws.Range("C25:D25").Formula = "=C12+K12"

# --- Row 26: yearly totals (MMscf/day * 365 / 1e6 style roll-up) ---
$ws.Range("C26").Formula = "=SUM(C16:C25)*365/1000000"
$ws.Range("D26").Formula = "=SUM(D16:D25)*365/1000000"

$ws.Range("A16:D26").NumberFormat = "0"

[void]$ws.Range("L20").Select()
